$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shopping section (rows 30-35): fill in the static-analysis results ---
# Cell values typed in the same order the author entered them, so the
# shared-string table comes out in the same sequence.
$ws.Range("F31").Value = "has CertificatePinner.smali, has 2 cert file including root, default network security config,"
$ws.Range("F34").Value = "using only http"
$ws.Range("F30").Value = "has CertificatePinner.smali, has 10 sha256 on networking/l.smali same file containing of okhttp, cleartext not permitted for dimains except system"
$ws.Range("F33").Value = "has CertificatePinner.smali and CertificatePinner.kt"
$ws.Range("F32").Value = 'has 2 cert file includiing root, has "Certificate pinning failure!" string, default security config, has all possible TLS handshake combination as string'
$ws.Range("F35").Value = 'defalt cecurity config, has string "Certificate pinning failure!", has cert factory and 9 end cert'

# Row 30 - 1. Ebay
$ws.Range("C30").Value = 1
# Row 31 - 2. CanadianTire
$ws.Range("C31").Value = 1
# Row 32 - 3. Walmart
$ws.Range("C32").Value = 1
# Row 33 - 4. Kijiji
$ws.Range("C33").Value = 1
# Row 34 - 5. Easyponno Seller
$ws.Range("C34").Value = 0
# Row 35 - 6. Food Panda
$ws.Range("C35").Value = 1

# --- Cosmetic cleanup on B25 (no longer highlighted) ---
$ws.Range("B25").Style = "Normal"

# --- View / selection state ---
$ws.Range("F4").Select() | Out-Null
